# Add season-record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108
$xlTop = -4160
$xlContinuous = 1
$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10

# New header cells, styled to match the existing header row (bold, thin box
# border, centered horizontal / top vertical alignment).
$newHeaders = @("AD1", "AE1", "AF1")
$newHeaderText = @("Wins", "Losses", "Ties")

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $cell = $ws.Range($newHeaders[$i])
    $cell.Value = $newHeaderText[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = $xlCenter
    $cell.VerticalAlignment = $xlTop
    $cell.Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
    $cell.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
    $cell.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
    $cell.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
}

# Every player row gets the team's season record: 65 wins, 97 losses, 0 ties.
for ($row = 2; $row -le 54; $row++) {
    $ws.Cells.Item($row, 30).Value = 65  # AD - Wins
    $ws.Cells.Item($row, 31).Value = 97  # AE - Losses
    $ws.Cells.Item($row, 32).Value = 0   # AF - Ties
}
